$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "H 72" record (row 2) was removed from the dataset; deleting the
# entire row shifts every subsequent row up by one, matching the new
# A1:F62 extent.
$ws.Rows(2).Delete()
